# Weekly refresh of the "Hortaliza, Vega Modelo de Temuco - Espinaca" data:
# a new daily record is inserted at row 259 (pushing the existing rows
# 259-292 down to 260-293), and the new row is populated with the latest
# reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 259 - this shifts rows 259:292
# down to 260:293 and bumps the sheet's used range to A1:R293.
$ws.Rows.Item(259).Insert()

# Populate the newly inserted row 259 with the new observation.
$ws.Cells.Item(259, 1).Value = 10
$ws.Cells.Item(259, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(259, 3).Value = "La Araucanía"
$ws.Cells.Item(259, 4).Value = 45127
$ws.Cells.Item(259, 5).Value = 9
$ws.Cells.Item(259, 6).Value = 100112012
$ws.Cells.Item(259, 7).Value = "Espinaca"
$ws.Cells.Item(259, 8).Value = "Sin especificar"
$ws.Cells.Item(259, 9).Value = "Primera"
$ws.Cells.Item(259, 10).Value = 30
$ws.Cells.Item(259, 11).Value = 8000
$ws.Cells.Item(259, 12).Value = 8000
$ws.Cells.Item(259, 13).Value = 8000
$ws.Cells.Item(259, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(259, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(259, 16).Value = 667
$ws.Cells.Item(259, 17).Value = 12
$ws.Cells.Item(259, 18).Value = "Hortaliza"
